$wb = $excel.ActiveWorkbook
$fitting = $wb.Worksheets.Item("Fitting")
$inference = $wb.Worksheets.Item("Inference")

# --- Fitting sheet: replace the questionnaire rows (2-18) with the new
# set of questions/answers (dataset/predictor/target/preprocessing
# questions replacing the old fitting-metric/model/algorithm questions). ---
$fitting.Cells.Item(2, 2).Value = "Dataset"
$fitting.Cells.Item(2, 3).Value = "Is the computational processing of the fitting dataset equivalent?"
$fitting.Cells.Item(2, 4).Value = "Yes. The method of preprocessing the data is the same and is carried out in a fully equivalent computational environment"
$fitting.Cells.Item(2, 5).Value = "Almost. The method includes minor variations in implementation, resulting in approximately similar outcomes"
$fitting.Cells.Item(2, 6).Value = "Somewhat. The method for preprocessing the data includes moderate differences in implementation, leading to materially different outcomes"
$fitting.Cells.Item(2, 7).Value = "Not equivalent or applicable. The preprocessing is inculdes substantial differences or completely different computational approaches are used"

$fitting.Cells.Item(3, 2).Value = "Dataset"
$fitting.Cells.Item(3, 3).Value = "Are the predictors selected for model fitting equivalent?"
$fitting.Cells.Item(3, 4).Value = "Yes. The same set of predictors or  a fully equivalent set of predictors is used"
$fitting.Cells.Item(3, 5).Value = "Almost. There is a slight variation in the set of predictors used, but approximately similar outcomes are expected."
$fitting.Cells.Item(3, 6).Value = "Somewhat. The predictor set exhibits moderate differences and a material effect on the outcomes is expected"
$fitting.Cells.Item(3, 7).Value = "Not equivalent or applicable. A substantially different sets of predictors is utilized"

$fitting.Cells.Item(4, 2).Value = "Dataset"
$fitting.Cells.Item(4, 3).Value = "Is the method used to select the predictors used for fitting equivalent?"
$fitting.Cells.Item(4, 4).Value = "Yes. The same strategy or a fully equivalent stratefy for selecting the predictors is used"
$fitting.Cells.Item(4, 5).Value = "Almost. There is slight variation in the strategy used to select the predictors, but approximately similar outcomes are expected"
$fitting.Cells.Item(4, 6).Value = "Somewhat. The strategy used to select the predictors exhibits moderate differences and a material difference in the outcomes is expected"
$fitting.Cells.Item(4, 7).Value = "Not equivalent or applicable. The strategy used to select predictors is substantially different"

$fitting.Cells.Item(5, 2).Value = "Dataset"
$fitting.Cells.Item(5, 3).Value = "Is the definition of the target event equivalent?"
$fitting.Cells.Item(5, 4).Value = "Yes. The target event definitions are identical or fully equivalent"
$fitting.Cells.Item(5, 5).Value = "Almost. There is slight variation in the target event definitions, but approximately similar outcomes are expected"
$fitting.Cells.Item(5, 6).Value = "Somewhat. There are moderate differences in target event definitions, and a material effect on the outcomes is expected"
$fitting.Cells.Item(5, 7).Value = "Not equivalent or applicable. The target definition is substantially different."

$fitting.Cells.Item(6, 2).Value = "Dataset"
$fitting.Cells.Item(6, 3).Value = "Are the preprocessing steps applied to the predictors before fitting equivalent?"
$fitting.Cells.Item(6, 4).Value = "Yes.  The same or  fully equivalent set of preprocessing steps is used"
$fitting.Cells.Item(6, 5).Value = "Almost.  There is a slight variation in the preprocessing steps being used, but approximately similar outcomes are expected"
$fitting.Cells.Item(6, 6).Value = "Somewhat.There are moderate differences in the data preprocessing steps, and a material effect on the outcomes is expected"
$fitting.Cells.Item(6, 7).Value = "Not equivalent or applicable."

$fitting.Cells.Item(7, 2).Value = "Dataset"
$fitting.Cells.Item(7, 3).Value = "Are the preprocessing steps applied to the target variable(s) before fitting equivalent?"
$fitting.Cells.Item(7, 4).Value = "Yes."
$fitting.Cells.Item(7, 5).Value = "Almost."
$fitting.Cells.Item(7, 6).Value = "Somewhat."
$fitting.Cells.Item(7, 7).Value = "Not equivalent or applicable."

$fitting.Cells.Item(8, 2).Value = "Dataset"
$fitting.Cells.Item(8, 3).Value = "Is the dataset size used for model fitting equivalent?"
$fitting.Cells.Item(8, 4).Value = "Yes."
$fitting.Cells.Item(8, 5).Value = "Almost."
$fitting.Cells.Item(8, 6).Value = "Somewhat."
$fitting.Cells.Item(8, 7).Value = "Not equivalent or applicable."

$fitting.Cells.Item(9, 2).Value = "Dataset"
$fitting.Cells.Item(9, 3).Value = "Is the method for partitioning the dataset into folds for use in fitting equivalent?"
$fitting.Cells.Item(9, 4).Value = "Yes."
$fitting.Cells.Item(9, 5).Value = "Almost."
$fitting.Cells.Item(9, 6).Value = "Somewhat."
$fitting.Cells.Item(9, 7).Value = "Not equivalent or applicable."

$fitting.Cells.Item(10, 2).Value = "Metrics"
$fitting.Cells.Item(10, 3).Value = "Is the way the evaluation metrics are calculated or implemented for model fitting equivalent?"
$fitting.Cells.Item(10, 4).Value = "Yes."
$fitting.Cells.Item(10, 5).Value = "Almost."
$fitting.Cells.Item(10, 6).Value = "Somewhat."
$fitting.Cells.Item(10, 7).Value = "Not equivalent or applicable."

$fitting.Cells.Item(11, 2).Value = "Metrics"
$fitting.Cells.Item(11, 3).Value = "Is the set of metrics chosen for use in model fitting equivalent?"
$fitting.Cells.Item(11, 4).Value = "Yes."
$fitting.Cells.Item(11, 5).Value = "Almost."
$fitting.Cells.Item(11, 6).Value = "Somewhat."
$fitting.Cells.Item(11, 7).Value = "Not equivalent or applicable."

$fitting.Cells.Item(12, 2).Value = "Metrics"
$fitting.Cells.Item(12, 3).Value = "Is the use of the chosen metrics within fitting optimization equivalent?"
$fitting.Cells.Item(12, 4).Value = "Yes."
$fitting.Cells.Item(12, 5).Value = "Almost."
$fitting.Cells.Item(12, 6).Value = "Somewhat."
$fitting.Cells.Item(12, 7).Value = "Not equivalent or applicable."

$fitting.Cells.Item(13, 2).Value = "Metrics"
$fitting.Cells.Item(13, 3).Value = "Are the observed va;ues of the metrics chosen for fitting equivalent?"
$fitting.Cells.Item(13, 4).Value = "Yes."
$fitting.Cells.Item(13, 5).Value = "Almost."
$fitting.Cells.Item(13, 6).Value = "Somewhat."
$fitting.Cells.Item(13, 7).Value = "Not equivalent or applicable."

$fitting.Cells.Item(14, 2).Value = "Model"
$fitting.Cells.Item(14, 3).Value = "Is the computational implementation of the algorithm equivalent?"
$fitting.Cells.Item(14, 4).Value = "Yes."
$fitting.Cells.Item(14, 5).Value = "Almost."
$fitting.Cells.Item(14, 6).Value = "Somewhat."
$fitting.Cells.Item(14, 7).Value = "Not equivalent or applicable."

$fitting.Cells.Item(15, 2).Value = "Model"
$fitting.Cells.Item(15, 3).Value = "Is the range of hyperparameters explored equivalent?"
$fitting.Cells.Item(15, 4).Value = "Yes."
$fitting.Cells.Item(15, 5).Value = "Almost."
$fitting.Cells.Item(15, 6).Value = "Somewhat."
$fitting.Cells.Item(15, 7).Value = "Not equivalent or applicable."

$fitting.Cells.Item(16, 2).Value = "Model"
$fitting.Cells.Item(16, 3).Value = "Is the stragegy used to search the hyperparameter space equivalent?"
$fitting.Cells.Item(16, 4).Value = "Yes."
$fitting.Cells.Item(16, 5).Value = "Almost."
$fitting.Cells.Item(16, 6).Value = "Somewhat."
$fitting.Cells.Item(16, 7).Value = "Not equivalent or applicable."

$fitting.Cells.Item(17, 2).Value = "Model"
$fitting.Cells.Item(17, 3).Value = "Are the best hyperparameters discovered during the search equivalent?"
$fitting.Cells.Item(17, 4).Value = "Yes."
$fitting.Cells.Item(17, 5).Value = "Almost."
$fitting.Cells.Item(17, 6).Value = "Somewhat."
$fitting.Cells.Item(17, 7).Value = "Not equivalent or applicable."

$fitting.Cells.Item(18, 2).Value = "Model"
$fitting.Cells.Item(18, 3).Value = "Are the learnable parameters observed after fitting equivalent?"
$fitting.Cells.Item(18, 4).Value = "Yes."
$fitting.Cells.Item(18, 5).Value = "Almost."
$fitting.Cells.Item(18, 6).Value = "Somewhat."
$fitting.Cells.Item(18, 7).Value = "Not equivalent or applicable."

# Rows 19-22 previously held questions 18-21; they are now blank (matching
# the already-empty rows below them), so clear out columns A:G.
$fitting.Range("A19:G22").ClearContents()

# Restore the Fitting sheet's remembered selection/active cell to match the
# post-edit state (also keeps "Fitting" as the active/tabSelected sheet).
$fitting.Range("L9").Select()
